# Update the crypto price/volume table with the latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.198.72'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.912.67'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.28'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5059'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3927'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09355'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.143'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.410'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.02'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.921.06'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.70%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001127'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06622'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.02'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9999'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.238'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.259.32'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.324'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.596'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.135.40'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.17'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '158.22'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.107'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.63%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.673'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.608'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.694'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06719'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02442'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2220'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.248'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.280'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6558'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.57'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.027'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9999'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6143'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.47'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.305'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.031'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '122.45'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.42%  '
